$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.497.69"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "2.311.54"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "105.88"
$ws.Range("E5").Value = "  +9.23%  "
$ws.Range("D6").Value = "308.73"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  +5.97%  "
$ws.Range("D10").Value = "36.29"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").Value = "52.95"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").Value = "0.0813"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "6.98"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").Value = "2.668.85"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").Value = "15.26"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").Value = "2.322.37"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "0.801"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "43.428.53"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").Value = "11.93"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "68.09"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "240.99"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("E25").Value = "  +4.57%  "
$ws.Range("D26").Value = "2.61"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "25.03"
$ws.Range("E28").Value = "  +7.24%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +5.52%  "
$ws.Range("D30").Value = "36.59"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.60"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").Value = "163.35"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D35").Value = "18.22"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  +6.42%  "
$ws.Range("D37").Value = "0.0735"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  +13.16%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  +4.28%  "
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "2.49"
$ws.Range("E43").Value = "  +14.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("D45").Value = "1.961.21"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").Value = "18.91"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").Value = "3.06"
$ws.Range("E47").Value = "  +5.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.30"
$ws.Range("E48").Value = "  +6.62%  "
$ws.Range("D49").Value = "58.07"
$ws.Range("E49").Value = "  +7.21%  "
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").Value = "1.58"
$ws.Range("E51").Value = "  +7.34%  "
